# Correções diversas de problemas encontrados durante a aula.
#
# 1) The "Espaço Reservado para Data" (datetime) placeholder on the
#    slide master, every custom layout, and the notes master had its
#    cached field text "05/05/2012" updated to "30/6/2012".
# 2) On slide 11, the title placeholder ("Classe interna anônima
#    (exemplo)") was narrowed (new width 8258204 EMU, was 8507288 EMU)
#    and its run was bumped to a 40pt (sz="4000") font.

$p = $ppt.ActivePresentation

$oldDate = "05/05/2012"
$newDate = "30/6/2012"

function Set-PlaceholderDate($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder.
Set-PlaceholderDate $p.SlideMaster

# Every custom layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Set-PlaceholderDate $layouts.Item($L)
}

# Notes master date placeholder: its Shapes collection does not allow
# in-place text edits in this host, but the HeadersFooters facade does.
$p.NotesMaster.HeadersFooters.DateAndTime.Text = $newDate

# Slide 11 ("Classe interna anônima (exemplo)") title shape tweaks.
$slide = $p.Slides.Item(11)
$title = $slide.Shapes.Item(2)
if ($title.TextFrame.TextRange.Text -eq "Classe interna anônima (exemplo)") {
    # Convert the new EMU width to points (1 pt = 12700 EMU); nudge by
    # half an EMU so double->float rounding lands back on the exact
    # integer EMU value instead of truncating one unit short.
    $title.Width = (8258204 + 0.5) / 12700.0
    $title.TextFrame.TextRange.Font.Size = 40
}
